$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry (2021-10-07) was added for this market/product.
# It lands at row 33, pushing the previously-existing rows 33-38 down to 34-39.
$ws.Rows("33:33").Insert()

$ws.Range("A33").Value = 2
$ws.Range("B33").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 44476
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100112026
$ws.Range("G33").Value = "Haba"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 1100
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = 5500
$ws.Range("N33").Value = "$/saco 25 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 220
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
